# Auto-generated edit script: applies scheduled market-data refresh
# to the per-job Leve profit sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW).
# For every touched row we rewrite the price/profit columns (H-N)
# with the refreshed values; a couple of rows also gain or lose
# their LeveProfitHQ (N) cell entirely, matching the upstream diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1400.5
$ws.Range("I40").Value = 1267.3334
$ws.Range("K40").Value = 1267.3334
$ws.Range("M40").Value = -1092.3334

$ws.Range("H80").Value = 2320.0881
$ws.Range("I80").Value = 736.94446
$ws.Range("J80").Value = 4101.125
$ws.Range("K80").Value = 2210.83338
$ws.Range("L80").Value = 12303.375
$ws.Range("M80").Value = -1212.83338
$ws.Range("N80").Value = -14299.375

$ws.Range("H82").Value = 6116.6665
$ws.Range("I82").Value = 3466.6667
$ws.Range("J82").Value = 8766.666999999999
$ws.Range("K82").Value = 10400.0001
$ws.Range("L82").Value = 26300.001
$ws.Range("M82").Value = -9994.000100000001
$ws.Range("N82").Value = -27112.001

$ws.Range("H83").Value = 2320.0881
$ws.Range("I83").Value = 736.94446
$ws.Range("J83").Value = 4101.125
$ws.Range("K83").Value = 6632.50014
$ws.Range("L83").Value = 36910.125
$ws.Range("M83").Value = -1640.50014
$ws.Range("N83").Value = -46894.125

$ws.Range("H85").Value = 6116.6665
$ws.Range("I85").Value = 3466.6667
$ws.Range("J85").Value = 8766.666999999999
$ws.Range("K85").Value = 10400.0001
$ws.Range("L85").Value = 26300.001
$ws.Range("M85").Value = -8996.000100000001
$ws.Range("N85").Value = -29108.001

$ws.Range("H100").Value = 12347392
$ws.Range("I100").Value = 16667560
$ws.Range("J100").Value = 4055.1428
$ws.Range("K100").Value = 16667560
$ws.Range("L100").Value = 4055.1428
$ws.Range("M100").Value = -16667019
$ws.Range("N100").Value = -5137.1428

$ws.Range("H106").Value = 2578.077
$ws.Range("I106").Value = 2201.6667
$ws.Range("J106").Value = 3425
$ws.Range("K106").Value = 2201.6667
$ws.Range("L106").Value = 3425
$ws.Range("M106").Value = -1570.6667
$ws.Range("N106").Value = -4687

$ws.Range("H137").Value = 1938.5238
$ws.Range("I137").Value = 1967.1666
$ws.Range("J137").Value = 1766.6666
$ws.Range("K137").Value = 5901.4998
$ws.Range("L137").Value = 5299.9998
$ws.Range("M137").Value = -3351.4998
$ws.Range("N137").Value = -10399.9998

$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 384046.16
$ws.Range("I32").Value = 2640.4854
$ws.Range("J32").Value = 4089129.8
$ws.Range("K32").Value = 2640.4854
$ws.Range("L32").Value = 4089129.8
$ws.Range("M32").Value = -2353.4854
$ws.Range("N32").Value = -4089703.8

$ws.Range("H74").Value = 610.2692
$ws.Range("I74").Value = 567.1053000000001
$ws.Range("K74").Value = 567.1053000000001
$ws.Range("M74").Value = 306.8946999999999

$ws.Range("H77").Value = 610.2692
$ws.Range("I77").Value = 567.1053000000001
$ws.Range("K77").Value = 2835.5265
$ws.Range("M77").Value = 1532.4735

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2631.125
$ws.Range("J20").Value = 2930.3
$ws.Range("L20").Value = 2930.3
$ws.Range("N20").Value = -3424.3

$ws.Range("H94").Value = 1474.2858
$ws.Range("I94").Value = 1561.8
$ws.Range("J94").Value = 1255.5
$ws.Range("K94").Value = 1561.8
$ws.Range("L94").Value = 1255.5
$ws.Range("M94").Value = -1110.8
$ws.Range("N94").Value = -2157.5

$ws.Range("H140").Value = 72695
$ws.Range("J140").Value = 72695
$ws.Range("L140").Value = 72695
$ws.Range("N140").Value = -83055

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2319.8918
$ws.Range("I31").Value = 1707.9286
$ws.Range("J31").Value = 4223.778
$ws.Range("K31").Value = 1707.9286
$ws.Range("L31").Value = 4223.778
$ws.Range("M31").Value = -1412.9286
$ws.Range("N31").Value = -4813.778

$ws.Range("H34").Value = 2319.8918
$ws.Range("I34").Value = 1707.9286
$ws.Range("J34").Value = 4223.778
$ws.Range("K34").Value = 1707.9286
$ws.Range("L34").Value = 4223.778
$ws.Range("M34").Value = -1505.9286
$ws.Range("N34").Value = -4627.778

$ws.Range("H58").Value = 797.43335
$ws.Range("I58").Value = 730.375
$ws.Range("K58").Value = 730.375
$ws.Range("M58").Value = -527.375

$ws.Range("H105").Value = 925.36365
$ws.Range("I105").Value = 953.2222
$ws.Range("K105").Value = 953.2222
$ws.Range("M105").Value = 793.7778

$ws.Range("H136").Value = 797.43335
$ws.Range("I136").Value = 730.375
$ws.Range("K136").Value = 2191.125
$ws.Range("M136").Value = 358.875

$ws.Range("H140").Value = 53833.332
$ws.Range("J140").Value = 53833.332
$ws.Range("L140").Value = 53833.332
$ws.Range("N140").Value = -64193.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1341.2
$ws.Range("I5").Value = 1341.2
$ws.Range("K5").Value = 4023.6
$ws.Range("M5").Value = -3911.6

$ws.Range("H74").Value = 13002.6
$ws.Range("I74").Value = 5013
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 15039
$ws.Range("L74").Value = 45000
$ws.Range("M74").Value = -13978
$ws.Range("N74").Value = -47122

$ws.Range("H77").Value = 13002.6
$ws.Range("I77").Value = 5013
$ws.Range("J77").Value = 15000
$ws.Range("K77").Value = 45117
$ws.Range("L77").Value = 135000
$ws.Range("M77").Value = -39813
$ws.Range("N77").Value = -145608

$ws.Range("H107").Value = 429.53333
$ws.Range("I107").Value = 329.83334
$ws.Range("K107").Value = 989.5000200000001
$ws.Range("M107").Value = 930.4999799999999

$ws.Range("H122").Value = 566
$ws.Range("I122").Value = 204
$ws.Range("J122").Value = 617.7143
$ws.Range("K122").Value = 1836
$ws.Range("L122").Value = 5559.428699999999
$ws.Range("M122").Value = 614
$ws.Range("N122").Value = -10459.4287

$ws.Range("H131").Value = 100003300
$ws.Range("J131").Value = 250005000
$ws.Range("L131").Value = 750015000
$ws.Range("N131").Value = -750025080

$ws.Range("H135").Value = 1341.2
$ws.Range("I135").Value = 1341.2
$ws.Range("K135").Value = 12070.8
$ws.Range("M135").Value = -9535.800000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17052858
$ws.Range("I70").Value = 40186036
$ws.Range("J70").Value = 7358.3687
$ws.Range("K70").Value = 40186036
$ws.Range("L70").Value = 7358.3687
$ws.Range("M70").Value = -40185766
$ws.Range("N70").Value = -7898.3687

$ws.Range("H73").Value = 17052858
$ws.Range("I73").Value = 40186036
$ws.Range("J73").Value = 7358.3687
$ws.Range("K73").Value = 40186036
$ws.Range("L73").Value = 7358.3687
$ws.Range("M73").Value = -40185100
$ws.Range("N73").Value = -9230.368699999999

$ws.Range("H80").Value = 4665.4
$ws.Range("I80").Value = 4330.8
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 4330.8
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -3332.8
$ws.Range("N80").Value = -6996

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H83").Value = 4665.4
$ws.Range("I83").Value = 4330.8
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 21654
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -16662
$ws.Range("N83").Value = -34984

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H102").Value = 1322.9667
$ws.Range("I102").Value = 1460.1052
$ws.Range("J102").Value = 1086.091
$ws.Range("K102").Value = 1460.1052
$ws.Range("L102").Value = 1086.091
$ws.Range("M102").Value = 161.8948
$ws.Range("N102").Value = -4330.091

$ws.Range("H126").Value = 13890869
$ws.Range("I126").Value = 1967.5
$ws.Range("K126").Value = 5902.5
$ws.Range("M126").Value = -3432.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2237.3547
$ws.Range("I7").Value = 1761.4166
$ws.Range("J7").Value = 3869.1428
$ws.Range("K7").Value = 1761.4166
$ws.Range("L7").Value = 3869.1428
$ws.Range("M7").Value = -1649.4166
$ws.Range("N7").Value = -4093.1428

$ws.Range("H126").Value = 2237.3547
$ws.Range("I126").Value = 1761.4166
$ws.Range("J126").Value = 3869.1428
$ws.Range("K126").Value = 5284.2498
$ws.Range("L126").Value = 11607.4284
$ws.Range("M126").Value = -2814.2498
$ws.Range("N126").Value = -16547.4284

$ws.Range("H132").Value = 20644.818
$ws.Range("I132").Value = 36693.434
$ws.Range("J132").Value = 1386.48
$ws.Range("K132").Value = 110080.302
$ws.Range("L132").Value = 4159.440000000001
$ws.Range("M132").Value = -107550.302
$ws.Range("N132").Value = -9219.440000000001

$ws.Range("H139").Value = 58120
$ws.Range("J139").Value = 58120
$ws.Range("L139").Value = 58120
$ws.Range("N139").Value = -68400

